$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update changed header / data cells ---
$ws.Range("O1").Value = "ABC1"
$ws.Range("P1").Value = "Abc12"
$ws.Range("Q1").Value = "Spcl Allowance"
$ws.Range("R1").Value = "Arrears"
$ws.Range("S1").Value = "Gross Pay"
$ws.Range("T1").Value = "PF"
$ws.Range("U1").Value = "ESIC"
$ws.Range("V1").Value = "PT"
$ws.Range("W1").Value = "TDS"
$ws.Range("X1").Value = "Deductible Arrears"
$ws.Range("Y1").Value = "total_deducations"
$ws.Range("Z1").Value = "NetPay"
$ws.Range("C2").Value = "balaraju Vankala"
$ws.Range("D2").Value = "27/11/2014"
$ws.Range("E2").Value = "New"
$ws.Range("H2").Value = 25000000.0
$ws.Range("I2").Value = 2083333.3333333333
$ws.Range("J2").Value = 1946880.0
$ws.Range("M2").Value = 833333.0
$ws.Range("N2").Value = 208333.0
$ws.Range("O2").Value = 100000.0
$ws.Range("P2").Value = 100000.0
$ws.Range("Q2").Value = 705208.0
$ws.Range("R2").Value = 0.0
$ws.Range("S2").Value = 1946880.0
$ws.Range("T2").Value = 100000.0
$ws.Range("U2").Value = 34070.3
$ws.Range("V2").Value = 0.0
$ws.Range("Y2").Value = 134070.0
$ws.Range("Z2").Value = 1812800.0
$ws.Range("C3").Value = "Balaraju vankala"
$ws.Range("E3").Value = "New"
$ws.Range("F3").Value = "HR Manager"
$ws.Range("G3").Value = "HR"
$ws.Range("H3").Value = 120000.0
$ws.Range("I3").Value = 10000.0
$ws.Range("J3").Value = 9345.0
$ws.Range("M3").Value = 4000.0
$ws.Range("N3").Value = 1000.0
$ws.Range("O3").Value = 480.0
$ws.Range("P3").Value = 480.0
$ws.Range("Q3").Value = 3385.0
$ws.Range("R3").Value = 0.0
$ws.Range("S3").Value = 9345.0
$ws.Range("T3").Value = 480.0
$ws.Range("U3").Value = 163.54
$ws.Range("V3").Value = 0.0
$ws.Range("Y3").Value = 643.54
$ws.Range("Z3").Value = 8701.46

# --- D3 needs special handling: keep literal text "2014-03-03" (not auto-parsed as a date) ---
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2014-03-03"
$ws.Range("D3").Style = "Normal"

# --- Remove trailing columns AA:AD (Deducted allowance1/2, old total_deducations/NetPay) ---
$ws.Range("AA1:AD3").EntireColumn.Delete()

# --- Re-apply column widths A:Z to match target layout ---
$ws.Columns.Item(1).ColumnWidth = 5.714285714285714
$ws.Columns.Item(2).ColumnWidth = 9.0
$ws.Columns.Item(3).ColumnWidth = 11.142857142857142
$ws.Columns.Item(4).ColumnWidth = 13.428571428571429
$ws.Columns.Item(5).ColumnWidth = 9.0
$ws.Columns.Item(6).ColumnWidth = 12.285714285714286
$ws.Columns.Item(7).ColumnWidth = 13.428571428571429
$ws.Columns.Item(8).ColumnWidth = 12.285714285714286
$ws.Columns.Item(9).ColumnWidth = 21.142857142857142
$ws.Columns.Item(10).ColumnWidth = 13.428571428571429
$ws.Columns.Item(11).ColumnWidth = 7.857142857142857
$ws.Columns.Item(12).ColumnWidth = 10.142857142857142
$ws.Columns.Item(13).ColumnWidth = 10.142857142857142
$ws.Columns.Item(14).ColumnWidth = 10.142857142857142
$ws.Columns.Item(15).ColumnWidth = 10.142857142857142
$ws.Columns.Item(16).ColumnWidth = 10.142857142857142
$ws.Columns.Item(17).ColumnWidth = 11.142857142857142
$ws.Columns.Item(18).ColumnWidth = 4.571428571428571
$ws.Columns.Item(19).ColumnWidth = 11.142857142857142
$ws.Columns.Item(20).ColumnWidth = 10.142857142857142
$ws.Columns.Item(21).ColumnWidth = 9.0
$ws.Columns.Item(22).ColumnWidth = 4.571428571428571
$ws.Columns.Item(23).ColumnWidth = 5.714285714285714
$ws.Columns.Item(24).ColumnWidth = 12.285714285714286
$ws.Columns.Item(25).ColumnWidth = 11.142857142857142
$ws.Columns.Item(26).ColumnWidth = 11.142857142857142
